$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update with latest event data
$ws.Range("A2:G2").NumberFormat = "@"
$ws.Range("K2:L2").NumberFormat = "@"
$ws.Range("A2").Value = '281474991152736-1748284058222'
$ws.Range("B2").Value = 'Mobile Usage'
$ws.Range("C2").Value = '2025-05-26T12:27:38.222'
$ws.Range("D2").Value = '281474991152736'
$ws.Range("E2").Value = '147'
$ws.Range("F2").Value = '51834001'
$ws.Range("G2").Value = 'MARCOS BARBOSA'
$ws.Range("H2").Value = 20.84346293
$ws.Range("I2").Value = -103.746949479
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1748284055722/p2gySLFM9k-camera-video-segment-driver-1748284058222.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSM7NHXLEK%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T140028Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIHkTN%2FfzOhHiomUjheOb%2F0CxYo9b4meFQT7IKd8z2rrqAiEAgFzZEKCoWLCzEK73cHmwpdAiC4Sr5XpZ4CBAYa9gwmYq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDG7jMTYOiXFy40SIXyq6A0LyUz37RICyGLNUjgsvRtr%2B1whULrl4lfJs2o135msL4TiQ4KEvgjvOgix7DVluiqvqgTqHJf0TiLcPNQthBvwtwiSAVa%2F%2FaBEdbCThE41bnqF8RQksiBEM0tlLBB5IvvAGgOmJFQm4SkA7yl%2F0%2F1KOHmq3E5NToWOWscz9SdaAlSMQS8cf%2B28c5nzPUomyb3u6EBpdLQTPXr1dcpVZzZkJ5t67WZefwtnxGMR9yOkDtqwL55cANIke15ai%2F5BnKuB%2B%2BJpx90Pjj6J9oWdV3ruV4QJd0zF1Aq6NtsHHEVGD6NTsIj9htTfn36a8E4RH3JG4rBFxqf4NCsYtD9%2BKnk45zz7z7OoXglpcNgBUptqdWv8QWiU4h5XzaOUwRcWOeZ3SOpV83NvYv2WmxMBipELujHlUMYjpIZK0915C2qgY5RSGkN9xTaIUJrBrqtUTsG%2FNors1TnX8hNeIKDKAK1GCjjr1ff7TXH%2BgK90%2FgQnUDx2AU8SSC1KO%2FWFkJjalpVd4D%2FpZZsWD2uBiJnx17zawkgcR6WwuMINfaItMI7DLO0VeU4ftu3CWpIAI2z6lFSG5ffnIcpQSR3Qw1YTXwQY6pQFctZXpaqbm0VnT8MdtFTOv5U8QCw2Itps%2BW4PERnsQP5DbMPiHSDRz9RmBR2ZtXyzgfO2GGmoF%2FAVZWC9fubu4%2Fzj46h2VOa8ZDb6e0ruEZUZ8iJz0cSYEFEXIwWQCbuTnh76t4UJiMBYgx56ShsmGIatbUVkqmo5cBNmERG0Hfx%2F9rUp7gtK0sbqrDxNGxU9nJ2UD2jbg28WS5QevZTZoOX994BE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2022%3A00%3A28%20GMT&X-Amz-Signature=7d69f221cbd0bc6a16d401f276d64f7b43e088977fb5746027ede31813ff41c4'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1748284055722/JuEVYDqJVY-camera-video-segment-1748284058222.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSM7NHXLEK%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T140028Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIHkTN%2FfzOhHiomUjheOb%2F0CxYo9b4meFQT7IKd8z2rrqAiEAgFzZEKCoWLCzEK73cHmwpdAiC4Sr5XpZ4CBAYa9gwmYq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDG7jMTYOiXFy40SIXyq6A0LyUz37RICyGLNUjgsvRtr%2B1whULrl4lfJs2o135msL4TiQ4KEvgjvOgix7DVluiqvqgTqHJf0TiLcPNQthBvwtwiSAVa%2F%2FaBEdbCThE41bnqF8RQksiBEM0tlLBB5IvvAGgOmJFQm4SkA7yl%2F0%2F1KOHmq3E5NToWOWscz9SdaAlSMQS8cf%2B28c5nzPUomyb3u6EBpdLQTPXr1dcpVZzZkJ5t67WZefwtnxGMR9yOkDtqwL55cANIke15ai%2F5BnKuB%2B%2BJpx90Pjj6J9oWdV3ruV4QJd0zF1Aq6NtsHHEVGD6NTsIj9htTfn36a8E4RH3JG4rBFxqf4NCsYtD9%2BKnk45zz7z7OoXglpcNgBUptqdWv8QWiU4h5XzaOUwRcWOeZ3SOpV83NvYv2WmxMBipELujHlUMYjpIZK0915C2qgY5RSGkN9xTaIUJrBrqtUTsG%2FNors1TnX8hNeIKDKAK1GCjjr1ff7TXH%2BgK90%2FgQnUDx2AU8SSC1KO%2FWFkJjalpVd4D%2FpZZsWD2uBiJnx17zawkgcR6WwuMINfaItMI7DLO0VeU4ftu3CWpIAI2z6lFSG5ffnIcpQSR3Qw1YTXwQY6pQFctZXpaqbm0VnT8MdtFTOv5U8QCw2Itps%2BW4PERnsQP5DbMPiHSDRz9RmBR2ZtXyzgfO2GGmoF%2FAVZWC9fubu4%2Fzj46h2VOa8ZDb6e0ruEZUZ8iJz0cSYEFEXIwWQCbuTnh76t4UJiMBYgx56ShsmGIatbUVkqmo5cBNmERG0Hfx%2F9rUp7gtK0sbqrDxNGxU9nJ2UD2jbg28WS5QevZTZoOX994BE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2022%3A00%3A28%20GMT&X-Amz-Signature=ba5edf960786960230fcf9697afefb207d5b56208d9d26ebfe6e52f821bcacf9'

# Row 3: update with latest event data
$ws.Range("A3:G3").NumberFormat = "@"
$ws.Range("K3:L3").NumberFormat = "@"
$ws.Range("A3").Value = '281474991152736-1748283981963'
$ws.Range("B3").Value = 'Mobile Usage'
$ws.Range("C3").Value = '2025-05-26T12:26:21.963'
$ws.Range("D3").Value = '281474991152736'
$ws.Range("E3").Value = '147'
$ws.Range("F3").Value = '51834001'
$ws.Range("G3").Value = 'MARCOS BARBOSA'
$ws.Range("H3").Value = 20.83726332
$ws.Range("I3").Value = -103.73671235
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1748283979463/3bY9s5rKOY-camera-video-segment-driver-1748283981963.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSM7NHXLEK%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T140028Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIHkTN%2FfzOhHiomUjheOb%2F0CxYo9b4meFQT7IKd8z2rrqAiEAgFzZEKCoWLCzEK73cHmwpdAiC4Sr5XpZ4CBAYa9gwmYq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDG7jMTYOiXFy40SIXyq6A0LyUz37RICyGLNUjgsvRtr%2B1whULrl4lfJs2o135msL4TiQ4KEvgjvOgix7DVluiqvqgTqHJf0TiLcPNQthBvwtwiSAVa%2F%2FaBEdbCThE41bnqF8RQksiBEM0tlLBB5IvvAGgOmJFQm4SkA7yl%2F0%2F1KOHmq3E5NToWOWscz9SdaAlSMQS8cf%2B28c5nzPUomyb3u6EBpdLQTPXr1dcpVZzZkJ5t67WZefwtnxGMR9yOkDtqwL55cANIke15ai%2F5BnKuB%2B%2BJpx90Pjj6J9oWdV3ruV4QJd0zF1Aq6NtsHHEVGD6NTsIj9htTfn36a8E4RH3JG4rBFxqf4NCsYtD9%2BKnk45zz7z7OoXglpcNgBUptqdWv8QWiU4h5XzaOUwRcWOeZ3SOpV83NvYv2WmxMBipELujHlUMYjpIZK0915C2qgY5RSGkN9xTaIUJrBrqtUTsG%2FNors1TnX8hNeIKDKAK1GCjjr1ff7TXH%2BgK90%2FgQnUDx2AU8SSC1KO%2FWFkJjalpVd4D%2FpZZsWD2uBiJnx17zawkgcR6WwuMINfaItMI7DLO0VeU4ftu3CWpIAI2z6lFSG5ffnIcpQSR3Qw1YTXwQY6pQFctZXpaqbm0VnT8MdtFTOv5U8QCw2Itps%2BW4PERnsQP5DbMPiHSDRz9RmBR2ZtXyzgfO2GGmoF%2FAVZWC9fubu4%2Fzj46h2VOa8ZDb6e0ruEZUZ8iJz0cSYEFEXIwWQCbuTnh76t4UJiMBYgx56ShsmGIatbUVkqmo5cBNmERG0Hfx%2F9rUp7gtK0sbqrDxNGxU9nJ2UD2jbg28WS5QevZTZoOX994BE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2022%3A00%3A28%20GMT&X-Amz-Signature=1fd3e3c2462906c865dd95a39616cea10111d88a0e630bf00e448c348e35a762'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1748283979463/qChm84icV1-camera-video-segment-1748283981963.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSM7NHXLEK%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T140028Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIHkTN%2FfzOhHiomUjheOb%2F0CxYo9b4meFQT7IKd8z2rrqAiEAgFzZEKCoWLCzEK73cHmwpdAiC4Sr5XpZ4CBAYa9gwmYq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDG7jMTYOiXFy40SIXyq6A0LyUz37RICyGLNUjgsvRtr%2B1whULrl4lfJs2o135msL4TiQ4KEvgjvOgix7DVluiqvqgTqHJf0TiLcPNQthBvwtwiSAVa%2F%2FaBEdbCThE41bnqF8RQksiBEM0tlLBB5IvvAGgOmJFQm4SkA7yl%2F0%2F1KOHmq3E5NToWOWscz9SdaAlSMQS8cf%2B28c5nzPUomyb3u6EBpdLQTPXr1dcpVZzZkJ5t67WZefwtnxGMR9yOkDtqwL55cANIke15ai%2F5BnKuB%2B%2BJpx90Pjj6J9oWdV3ruV4QJd0zF1Aq6NtsHHEVGD6NTsIj9htTfn36a8E4RH3JG4rBFxqf4NCsYtD9%2BKnk45zz7z7OoXglpcNgBUptqdWv8QWiU4h5XzaOUwRcWOeZ3SOpV83NvYv2WmxMBipELujHlUMYjpIZK0915C2qgY5RSGkN9xTaIUJrBrqtUTsG%2FNors1TnX8hNeIKDKAK1GCjjr1ff7TXH%2BgK90%2FgQnUDx2AU8SSC1KO%2FWFkJjalpVd4D%2FpZZsWD2uBiJnx17zawkgcR6WwuMINfaItMI7DLO0VeU4ftu3CWpIAI2z6lFSG5ffnIcpQSR3Qw1YTXwQY6pQFctZXpaqbm0VnT8MdtFTOv5U8QCw2Itps%2BW4PERnsQP5DbMPiHSDRz9RmBR2ZtXyzgfO2GGmoF%2FAVZWC9fubu4%2Fzj46h2VOa8ZDb6e0ruEZUZ8iJz0cSYEFEXIwWQCbuTnh76t4UJiMBYgx56ShsmGIatbUVkqmo5cBNmERG0Hfx%2F9rUp7gtK0sbqrDxNGxU9nJ2UD2jbg28WS5QevZTZoOX994BE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2022%3A00%3A28%20GMT&X-Amz-Signature=40153e66a104f5edf56ee4eed05c5d3f50cd0223bd194785d1f9dea0dc861e90'

# Row 4: update with latest event data
$ws.Range("A4:G4").NumberFormat = "@"
$ws.Range("K4:L4").NumberFormat = "@"
$ws.Range("A4").Value = '281474991152736-1748276160117'
$ws.Range("B4").Value = 'Mobile Usage'
$ws.Range("C4").Value = '2025-05-26T10:16:00.117'
$ws.Range("D4").Value = '281474991152736'
$ws.Range("E4").Value = '147'
$ws.Range("F4").Value = '51834001'
$ws.Range("G4").Value = 'MARCOS BARBOSA'
$ws.Range("H4").Value = 20.72414717
$ws.Range("I4").Value = -103.503679399
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1748276157617/6zxWnZYxXR-camera-video-segment-driver-1748276160117.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSM7NHXLEK%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T140029Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIHkTN%2FfzOhHiomUjheOb%2F0CxYo9b4meFQT7IKd8z2rrqAiEAgFzZEKCoWLCzEK73cHmwpdAiC4Sr5XpZ4CBAYa9gwmYq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDG7jMTYOiXFy40SIXyq6A0LyUz37RICyGLNUjgsvRtr%2B1whULrl4lfJs2o135msL4TiQ4KEvgjvOgix7DVluiqvqgTqHJf0TiLcPNQthBvwtwiSAVa%2F%2FaBEdbCThE41bnqF8RQksiBEM0tlLBB5IvvAGgOmJFQm4SkA7yl%2F0%2F1KOHmq3E5NToWOWscz9SdaAlSMQS8cf%2B28c5nzPUomyb3u6EBpdLQTPXr1dcpVZzZkJ5t67WZefwtnxGMR9yOkDtqwL55cANIke15ai%2F5BnKuB%2B%2BJpx90Pjj6J9oWdV3ruV4QJd0zF1Aq6NtsHHEVGD6NTsIj9htTfn36a8E4RH3JG4rBFxqf4NCsYtD9%2BKnk45zz7z7OoXglpcNgBUptqdWv8QWiU4h5XzaOUwRcWOeZ3SOpV83NvYv2WmxMBipELujHlUMYjpIZK0915C2qgY5RSGkN9xTaIUJrBrqtUTsG%2FNors1TnX8hNeIKDKAK1GCjjr1ff7TXH%2BgK90%2FgQnUDx2AU8SSC1KO%2FWFkJjalpVd4D%2FpZZsWD2uBiJnx17zawkgcR6WwuMINfaItMI7DLO0VeU4ftu3CWpIAI2z6lFSG5ffnIcpQSR3Qw1YTXwQY6pQFctZXpaqbm0VnT8MdtFTOv5U8QCw2Itps%2BW4PERnsQP5DbMPiHSDRz9RmBR2ZtXyzgfO2GGmoF%2FAVZWC9fubu4%2Fzj46h2VOa8ZDb6e0ruEZUZ8iJz0cSYEFEXIwWQCbuTnh76t4UJiMBYgx56ShsmGIatbUVkqmo5cBNmERG0Hfx%2F9rUp7gtK0sbqrDxNGxU9nJ2UD2jbg28WS5QevZTZoOX994BE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2022%3A00%3A29%20GMT&X-Amz-Signature=5f0a83de0e8d2d7f0190de20c52d438b0a07ebb1f8d3820bb2c05453fe35e963'
$ws.Range("L4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1748276157617/AqojMFPYnw-camera-video-segment-1748276160117.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSM7NHXLEK%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T140029Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIHkTN%2FfzOhHiomUjheOb%2F0CxYo9b4meFQT7IKd8z2rrqAiEAgFzZEKCoWLCzEK73cHmwpdAiC4Sr5XpZ4CBAYa9gwmYq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDG7jMTYOiXFy40SIXyq6A0LyUz37RICyGLNUjgsvRtr%2B1whULrl4lfJs2o135msL4TiQ4KEvgjvOgix7DVluiqvqgTqHJf0TiLcPNQthBvwtwiSAVa%2F%2FaBEdbCThE41bnqF8RQksiBEM0tlLBB5IvvAGgOmJFQm4SkA7yl%2F0%2F1KOHmq3E5NToWOWscz9SdaAlSMQS8cf%2B28c5nzPUomyb3u6EBpdLQTPXr1dcpVZzZkJ5t67WZefwtnxGMR9yOkDtqwL55cANIke15ai%2F5BnKuB%2B%2BJpx90Pjj6J9oWdV3ruV4QJd0zF1Aq6NtsHHEVGD6NTsIj9htTfn36a8E4RH3JG4rBFxqf4NCsYtD9%2BKnk45zz7z7OoXglpcNgBUptqdWv8QWiU4h5XzaOUwRcWOeZ3SOpV83NvYv2WmxMBipELujHlUMYjpIZK0915C2qgY5RSGkN9xTaIUJrBrqtUTsG%2FNors1TnX8hNeIKDKAK1GCjjr1ff7TXH%2BgK90%2FgQnUDx2AU8SSC1KO%2FWFkJjalpVd4D%2FpZZsWD2uBiJnx17zawkgcR6WwuMINfaItMI7DLO0VeU4ftu3CWpIAI2z6lFSG5ffnIcpQSR3Qw1YTXwQY6pQFctZXpaqbm0VnT8MdtFTOv5U8QCw2Itps%2BW4PERnsQP5DbMPiHSDRz9RmBR2ZtXyzgfO2GGmoF%2FAVZWC9fubu4%2Fzj46h2VOa8ZDb6e0ruEZUZ8iJz0cSYEFEXIwWQCbuTnh76t4UJiMBYgx56ShsmGIatbUVkqmo5cBNmERG0Hfx%2F9rUp7gtK0sbqrDxNGxU9nJ2UD2jbg28WS5QevZTZoOX994BE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2022%3A00%3A29%20GMT&X-Amz-Signature=29e9a120da0324c5861d0053824e31cb7481a4eec8292f21e6a64f3a7dd3dbdd'

# Remove the two oldest event rows (previously rows 5 and 6)
$ws.Rows(5).Delete()
$ws.Rows(5).Delete()

